# OBIE Products.xlsx - "Add files via upload"
#
# Row 7 (the AccountId row) is updated:
#   - H7: "Yes" -> "OB to populate"
#   - I7: the existing CustomerReference/PartyIdentification/IdentifierValue
#         mapping note gets a "Ref:" header line prepended.
# The active selection moves from I2 to I8, and row 7 grows to the same
# auto-fit height used by its wrapped-text neighbours (73.5pt) now that it
# holds a populated, wrapped I-column note like the surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "OB to populate"
$ws.Range("I7").Value = "Ref:`nCustomerReference:`n              `$ref: '#/components/schemas/Involvedparty'`nPartyIdentification`nIdentifierValue"

$ws.Rows.Item(7).RowHeight = 73.5

$ws.Range("I8").Select()
